# Add newly tracked apartment complexes to the bottom of the list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @(193, 8626,   "귀인마을현대홈타운"),
    @(194, 1483,   "초원2단지대림"),
    @(195, 1471,   "무궁화경남"),
    @(196, 107579, "평촌더샵센트럴시티"),
    @(197, 2505,   "인덕원마을삼성"),
    @(198, 1467,   "은하수(벽산)"),
    @(199, 154917, "평촌센텀퍼스트"),
    @(200, 142558, "평촌트리지아"),
    @(201, 3081,   "호계현대홈타운1차"),
    @(202, 126060, "평촌래미안푸르지오"),
    @(203, 102312, "동편마을3단지"),
    @(204, 144023, "평촌엘프라우드")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Scroll the view down to show the newly added rows and move the selection,
# matching where the user ended up after pasting the new data.
$excel.ActiveWindow.ScrollRow = 190
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D207").Select()
